$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

# Columns that must stay plain text even though they look numeric / date-like
# (Caso, F. De Reclamo, Comuna, OT) need NumberFormat forced to Text first so
# Excel doesn't silently coerce them into numbers / date serials.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "7516"

$ws.Range("B$row").NumberFormat = "@"
$ws.Range("B$row").Value = "10/15/2025"

$ws.Range("C$row").Value = "ALVAREZ, CRISOSTOMO 3000"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "7"

$ws.Range("E$row").NumberFormat = "@"
$ws.Range("E$row").Value = "810371027"

$ws.Range("F$row").Value = "Optical Power"
$ws.Range("G$row").Value = "Pendiente"
$ws.Range("H$row").Value = "Picada"

$ws.Range("I$row").Value = 1

$ws.Range("J$row").Value = "Cambio"
$ws.Range("K$row").Value = "Sin equipos"
$ws.Range("L$row").Value = "Pasante"

$ws.Range("M$row").Value = -58.458516
$ws.Range("N$row").Value = -34.646422

$ws.Range("O$row").Value = "Boedo"
$ws.Range("P$row").Value = "Capital Sur"
$ws.Range("Q$row").Value = "PPT-N"
$ws.Range("R$row").Value = "Fuera de Poligono OVL"
